# Issue no. 58 @1342808581
#
# Populate the "Umgesetzt" (implemented) / "Verworfen" (discarded) columns
# for a handful of rows on the "fare" sheet, then re-apply the existing
# autofilter (on columns G/I, "blank" only) by hiding the rows that now
# carry values. Finally move the active-cell selection to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fare")

$implementedDate = 41108   # 2012-07-18

# --- Row 74: "Projector, Epoch" -----------------------------------------
# Umgesetzt (implemented): date + "Epoch."
$ws.Range("G74").Value = $implementedDate

# Verworfen (discarded): date (custom "d-mmm" number format) + "Projector."
# (new shared strings are interned in the order they are first written, so
# "Projector." must land before "Epoch." to match the source order)
$ws.Range("J74").Value = "Projector."
$ws.Range("H74").Value = "Epoch."
$ws.Range("I74").NumberFormat = "d-mmm"
$ws.Range("I74").Value = $implementedDate

# --- Row 109 -------------------------------------------------------------
$ws.Range("G109").Value = $implementedDate

# --- Row 113 -------------------------------------------------------------
$ws.Range("G113").Value = $implementedDate

# --- Row 116 -------------------------------------------------------------
$ws.Range("G116").Value = $implementedDate
$ws.Range("H116").Value = "Interface Projector überarbeitet."

# Re-hide the rows that the G/I "blanks" autofilter would now exclude.
$ws.Rows.Item(74).Hidden = $true
$ws.Rows.Item(109).Hidden = $true
$ws.Rows.Item(113).Hidden = $true
$ws.Rows.Item(116).Hidden = $true

# Move the selection (bottom-right frozen pane) to G1.
[void]$ws.Activate()
[void]$ws.Range("G1").Select()
